$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")

$ws.Range("D2").Value = "48.098.73"
$ws.Range("E2").Value = "  +1.62%  "

$ws.Range("D3").Value = "2.505.71"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$helper.NumberFormat = "@"
$helper.Value = "320.57"
$helper.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$helper.NumberFormat = "@"
$helper.Value = "108.91"
$helper.Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +0.22%  "

$helper.NumberFormat = "@"
$helper.Value = "0.528"
$helper.Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  +1.21%  "

$ws.Range("E8").Value = "  +0.03%  "

$helper.NumberFormat = "@"
$helper.Value = "0.544"
$helper.Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +1.40%  "

$helper.NumberFormat = "@"
$helper.Value = "39.79"
$helper.Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +1.46%  "

$helper.NumberFormat = "@"
$helper.Value = "20.04"
$helper.Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +8.76%  "

$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("D15").Value = "2.901.47"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "2.517.10"
$ws.Range("E16").Value = "  +1.63%  "

$helper.NumberFormat = "@"
$helper.Value = "0.845"
$helper.Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").Value = "47.948.70"
$ws.Range("E18").Value = "  +1.50%  "

$ws.Range("E19").Value = "  +0.61%  "

$helper.NumberFormat = "@"
$helper.Value = "6.61"
$helper.Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("D21").Value = "0.0₃0942"
$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("E22").Value = "  +3.24%  "

$helper.NumberFormat = "@"
$helper.Value = "72.15"
$helper.Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +2.50%  "

$helper.NumberFormat = "@"
$helper.Value = "274.18"
$helper.Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +11.82%  "

$helper.NumberFormat = "@"
$helper.Value = "2.56"
$helper.Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("E26").Value = "  +0.02%  "

$helper.NumberFormat = "@"
$helper.Value = "25.86"
$helper.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("E28").Value = "  +5.57%  "

$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("E30").Value = "  +1.87%  "

$helper.NumberFormat = "@"
$helper.Value = "35.41"
$helper.Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +1.71%  "

$helper.NumberFormat = "@"
$helper.Value = "49.43"
$helper.Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -0.64%  "

$helper.NumberFormat = "@"
$helper.Value = "19.27"
$helper.Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -6.69%  "

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$helper.NumberFormat = "@"
$helper.Value = "1.00"
$helper.Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$helper.NumberFormat = "@"
$helper.Value = "5.34"
$helper.Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("E37").Value = "  -0.56%  "

$helper.NumberFormat = "@"
$helper.Value = "4.61"
$helper.Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -2.62%  "

$ws.Range("E39").Value = "  +0.95%  "

$ws.Range("E40").Value = "  +0.74%  "

$helper.NumberFormat = "@"
$helper.Value = "122.34"
$helper.Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +5.21%  "

$helper.NumberFormat = "@"
$helper.Value = "2.22"
$helper.Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -0.50%  "

$helper.NumberFormat = "@"
$helper.Value = "21.75"
$helper.Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -6.14%  "

$helper.NumberFormat = "@"
$helper.Value = "0.0305"
$helper.Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +3.07%  "

$ws.Range("D45").Value = "2.020.40"
$ws.Range("E45").Value = "  +1.38%  "

$helper.NumberFormat = "@"
$helper.Value = "3.12"
$helper.Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +2.67%  "

$ws.Range("E47").Value = "  +3.75%  "

$ws.Range("E48").Value = "  -0.63%  "

$helper.NumberFormat = "@"
$helper.Value = "9.02"
$helper.Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -1.49%  "

$ws.Range("E50").Value = "  +1.57%  "

$ws.Range("E51").Value = "  +2.53%  "

$helper.Clear() | Out-Null
$excel.CutCopyMode = 0